$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the OAMP model reference from OPA336UJ to OPA336N and its unit price.
$ws.Range("D7").Value = "OPA336N"
$ws.Range("E7").Value = 1.29

# Reflect the last user selection on the sheet.
$ws.Range("D7").Select()
